$d = $word.ActiveDocument

$replacements = @(
    @("27×67=", "40×46="),
    @("85×77=", "57×98="),
    @("41×21=", "97×42="),
    @("40×66=", "50×88="),
    @("95×75=", "27×31="),
    @("74×71=", "88×33="),
    @("44×53=", "45×56="),
    @("77×62=", "66×29="),
    @("35×45=", "52×70="),
    @("80×73=", "21×33="),
    @("76×15=", "23×16="),
    @("48×44=", "21×69="),
    @("29×75=", "65×19="),
    @("81×86=", "73×94="),
    @("24×59=", "73×85="),
    @("89×80=", "75×58="),
    @("45×64=", "85×53="),
    @("86×57=", "99×63="),
    @("20×92=", "99×52="),
    @("49×37=", "70×38="),
    @("52×92=", "34×43="),
    @("80×91=", "93×27="),
    @("79×52=", "89×27="),
    @("70×94=", "24×70="),
    @("29×55=", "63×57=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
